$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variavel = "Índice de Gini do rendimento domiciliar per capita, a preços médios do ano"
$ano = "31/12/2022"

# Full target data for rows 2..27 (A: Região, D: Valor, E: Colocação)
$data = @(
  @{A="Paraíba";             D=0.5580000000000001; E="1º"},
  @{A="Roraima";             D=0.547;               E="2º"},
  @{A="Rio de Janeiro";      D=0.54;                E="3º"},
  @{A="Distrito Federal";    D=0.536;               E="4º"},
  @{A="Amapá";               D=0.531;               E="5º"},
  @{A="Sergipe";             D=0.528;               E="6º"},
  @{A="Ceará";               D=0.518;               E="10º"},
  @{A="Pernambuco";          D=0.515;               E="11º"},
  @{A="Bahia";               D=0.511;               E="12º"},
  @{A="Amazonas";            D=0.509;               E="13º"},
  @{A="Pará";                D=0.508;               E="14º"},
  @{A="Tocantins";           D=0.507;               E="15º"},
  @{A="São Paulo";           D=0.5;                 E="16º"},
  @{A="Alagoas";             D=0.498;               E="17º"},
  @{A="Espírito Santo";      D=0.493;               E="18º"},
  @{A="Maranhão";            D=0.491;               E="19º"},
  @{A="Mato Grosso do Sul";  D=0.478;               E="20º"},
  @{A="Paraná";              D=0.47;                E="21º"},
  @{A="Rio Grande do Sul";   D=0.467;               E="22º"},
  @{A="Minas Gerais";        D=0.466;               E="23º"},
  @{A="Goiás";               D=0.456;               E="24º"},
  @{A="Mato Grosso";         D=0.45;                E="25º"},
  @{A="Rondônia";            D=0.447;               E="26º"},
  @{A="Santa Catarina";      D=0.419;               E="27º"},
  @{A="Nordeste";            D=0.517;               E=$null},
  @{A="Brasil";              D=0.518;               E=$null}
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $variavel
    $ws.Cells.Item($r, 3).Value2 = $ano
    $ws.Cells.Item($r, 4).Value2 = $row.D
    if ($row.E -ne $null) {
        $ws.Cells.Item($r, 5).Value2 = $row.E
    } else {
        $ws.Cells.Item($r, 5).ClearContents()
    }
    $r++
}
